# edit.ps1 - apply "Audit of Operational Petty Cash Management" content rewrite
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - Title slide
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Title
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Audit of Operational Petty Cash Management"

# Content placeholder: collapse "Report:/Business Unit:/Date:" block into two
# bold lines, each with extra space-before, and drop the old 3rd paragraph.
$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange

$tr1.Paragraphs(1).Runs(1).Text = "Regional Branch Operations"
$tr1.Paragraphs(1).Runs(2).Text = ""
$tr1.Paragraphs(1).ParagraphFormat.SpaceBefore = 30

$tr1.Paragraphs(2).Runs(1).Text = "March 20, 2026"
$tr1.Paragraphs(2).Runs(2).Text = ""
$tr1.Paragraphs(2).ParagraphFormat.SpaceBefore = 30

$tr1.Paragraphs(3).Delete()

# ---------------------------------------------------------------------------
# Slide 2 - Executive Summary
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange

$tr2.Paragraphs(2).Text = "To verify the security of petty cash funds and the validity of small-value disbursements."
$tr2.Paragraphs(4).Text = "15 regional branches maintain petty cash floats ranging from `$500 to `$2,000."
$tr2.Paragraphs(6).Text = "Surprise cash counts and voucher audits at 5 selected regional branches."

# ---------------------------------------------------------------------------
# Slide 3 - Observation 1 (was "Commingling of Personal Funds", now
# "Missing Reconciliation Logs")
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Observation 1: Missing Reconciliation Logs"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange

$tr3.Paragraphs(1).Runs(2).Text = " 3 out of 5 branches visited did not perform weekly cash reconciliations as required by policy."
$tr3.Paragraphs(2).Runs(2).Text = " Undetected theft or loss of funds; financial reporting inaccuracies."
# Paragraph 3 (Risk Rating: INADEQUATE) is unchanged.
$tr3.Paragraphs(4).Runs(2).Text = " Mandate a weekly sign-off by the Branch Manager on a standardized cash reconciliation form."
# Drop the trailing "Status: Open" paragraph entirely.
$tr3.Paragraphs(5).Delete()

# ---------------------------------------------------------------------------
# Slide 4 - Observation 2 (was "Missing Custodian Segregation", now
# "Commingling of Personal and Company Funds")
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Observation 2: Commingling of Personal and Company Funds"

$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange

$tr4.Paragraphs(1).Runs(2).Text = " Personal funds and ‘IOU’ notes from employees were found in the petty cash lockbox at one location."
$tr4.Paragraphs(2).Runs(2).Text = " High risk of fraud and lack of accountability for corporate assets."
$tr4.Paragraphs(3).Runs(3).Text = "INADEQUATE"
$tr4.Paragraphs(4).Runs(2).Text = " Strictly prohibit ‘IOUs’ and enforce immediate disciplinary action for commingling of funds."
# Drop the trailing "Status: In Progress" paragraph entirely.
$tr4.Paragraphs(5).Delete()

# ---------------------------------------------------------------------------
# Slide 5 - Recommendations Summary
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange

$tr5.Paragraphs(1).Text = "Standardize petty cash reconciliation templates."
$tr5.Paragraphs(2).Text = "Perform monthly unannounced cash counts by Regional Controllers."
$tr5.Paragraphs(3).Text = "Phase out physical petty cash in favor of corporate ‘P-Cards’ for small spend."

# ---------------------------------------------------------------------------
# Slide 6 - Management Action Plan (table)
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tbl = $s6.Shapes.Item(2).Table

# Insert a new row after the header for "P-Card Rollout".
$newRow = $tbl.Rows.Add(2)
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "P-Card Rollout"
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "Finance VP"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "September 2026"

# "Surprise Counts" row (now row 3): cadence Continuous -> Monthly.
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "Monthly"

# "Float Reduction" row (now row 4) becomes "Policy Briefing".
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Policy Briefing"
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "Branch Managers"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Text = "April 2026"

# Remove the old "Digital Logs" row (now row 5).
$tbl.Rows(5).Delete()
